$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs of the "中雨..." paragraph (paragraph 4) back into a
#    single run. Doing a Find/Replace across both runs (with the replacement
#    text identical to the concatenation of the two original runs) makes the
#    engine rebuild that span as one run and drops the now-redundant
#    "_GoBack" bookmark that used to sit between them.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "中雨，今天农历是五月初四，明天就是端午节了。", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "中雨，今天农历是五月初四，明天就是端午节了。", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Split the "多云，今天是六一儿童节，又是开心的一天呢" paragraph
#    (paragraph 2) into two runs, right after "...六一", and drop the
#    "_GoBack" bookmark there (Word only keeps a single "_GoBack" bookmark,
#    so adding a new one automatically removes the old one from paragraph 4).
# ---------------------------------------------------------------------------
$splitPoint = $d.Content
$splitPoint.Find.Execute(
    "多云，今天是六一", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$splitPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# ---------------------------------------------------------------------------
# 3) Paragraph 4 needs to grow into three paragraphs:
#      - the existing "中雨，今天农历是五月初四，明天就是端午节了。" text,
#        but with its paragraph-mark formatting hint switched from
#        "default" to "eastAsia"
#      - a new paragraph "2022年6月3日星期五"
#      - a new paragraph "中雨，今天是农历五月初五，中国传统端午节。"
#
#    Range.InsertXML() replaces the whole paragraph that the (collapsed)
#    range sits in, which is exactly what is needed to rewrite paragraph 4
#    in place as three paragraphs with precise pPr/rPr control (the
#    rFonts/@w:hint attribute isn't reachable through any Font/Paragraph
#    property, only through raw OOXML).
#
#    Because paragraph 4 is currently the last paragraph in the document,
#    InsertXML would also leave a stray empty trailing paragraph behind (the
#    body's implicit final mark). To avoid that, a temporary paragraph is
#    appended first so paragraph 4 is no longer last; that temporary
#    paragraph is removed again at the end.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$tempAnchor = $p4.Range
$tempAnchor.Collapse(0)
$tempAnchor.InsertParagraphAfter() | Out-Null

$p4 = $d.Paragraphs(4)
$target = $p4.Range
$target.Collapse(1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>中雨，今天农历是五月初四，明天就是端午节了。</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2022年6月3日星期五</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>中雨，今天是农历五月初五，中国传统端午节。</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml) | Out-Null

# Remove the temporary trailing empty paragraph (now the last one), together
# with the paragraph mark before it, so the paragraph count ends up correct.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$cleanup = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$cleanup.Delete() | Out-Null
